# Fix pre-analyser val test with MCDS 7.4
# This script reproduces the target diff:
#  - On "Details" sheet: insert 11 new columns (BO..BY) before the former
#    "EstA(1)" column, fill their header row (row 1) and their data for
#    rows 2..13, and update the existing T/U/V columns (run timestamp,
#    run duration, run folder name) for rows 2..13.
#  - On "Synthesis" sheet: update the AI column (RunFolder) for rows 2..13
#    to match the new run folder names.
#  - On "Computing platform" sheet: update the pyaudisam version and the
#    MCDS engine path (B10, B11).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Details" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Details")

# Insert 11 new columns before the old "BO" column (EstA(1)); this shifts
# the former BO:CL range to BZ:CW, matching dimension A1:CW13 afterwards.
$ws.Range("BO1:BY1").EntireColumn.Insert()

# New header labels for the 11 freshly inserted columns (row 1)
$newHeaders = @{
    "BO1" = "SansDoc #18"
    "BP1" = "SansDoc #19"
    "BQ1" = "SansDoc #20 1"
    "BR1" = "SansDoc #20 2"
    "BS1" = "SansDoc #20 3"
    "BT1" = "SansDoc #21 1"
    "BU1" = "SansDoc #21 2"
    "BV1" = "SansDoc #21 3"
    "BW1" = "SansDoc #22 1"
    "BX1" = "SansDoc #22 2"
    "BY1" = "SansDoc #22 3"
}
foreach ($addr in $newHeaders.Keys) {
    $ws.Range($addr).Value = $newHeaders[$addr]
}

# New data values for the 11 inserted columns, per data row (2..13), in
# column order BO,BP,BQ,BR,BS,BT,BU,BV,BW,BX,BY
$newColsData = @{
    2  = @(1527.339, 1, 1.264273, 1.566549, 1.322656, 4.812198, 9.740577, 12.20725, 3, 6, 11)
    3  = @(2340.707, 1, 1.654648, 2.34082, 1.79977, 9.927889, 23.96511, 28.82581, 6, 10, 16)
    4  = @(2872.073, 1, 2.937908, 2.978181, 2.835873, 18.90978, 35.17652, 53.57235, 6, 11, 18)
    5  = @(4350.866, 1, 1.462845, 1.885604, 1.785336, 14.17539, 29.57599, 46.39332, 9, 15, 24)
    6  = @(1581.477, 1, 0.5498409, 0.5520127, 1.207226, 2.267506, 2.852037, 13.09932, 4, 7, 11)
    7  = @(2684.004, 1, 2.053263, 1.781311, 1.267737, 3.047782, 6.986795, 12.02188, 5, 9, 15)
    8  = @(2745.426, 1, 0.6904538, 0.9984319, 1.003679, 2.713422, 8.238899999999999, 9.452434999999999, 6, 9, 14)
    9  = @(4797.853, 1, 3.137696, 2.30542, 1.867133, 6.597888, 12.92294, 23.27799, 8, 14, 21)
    10 = @(717.264, 1, 1.373853, 1.724357, 2.352862, 0.6097739, 1.368052, 6.376579, 1, 2, 5)
    11 = @(1054.446, 1, 0.6924571, 1.217724, 1.076941, 0.01215697, 2.662995, 1.371374, 2, 4, 7)
    12 = @(1333.428, 1, 1.121317, 1.361363, 1.588817, 1.354986, 2.235081, 5.116446, 2, 5, 8)
    13 = @(1944.133, 1, 0.8963301, 0.9087912, 1.087616, 2.099729, 2.986194, 6.712539, 4, 7, 11)
}
$newColsOrder = @("BO", "BP", "BQ", "BR", "BS", "BT", "BU", "BV", "BW", "BX", "BY")

foreach ($row in $newColsData.Keys) {
    $values = $newColsData[$row]
    for ($i = 0; $i -lt $newColsOrder.Length; $i++) {
        $addr = "$($newColsOrder[$i])$row"
        $ws.Range($addr).Value = $values[$i]
    }
}

# Updated run timestamp (T), run duration (U) and run folder name (V),
# for each data row 2..13
$tuvData = @{
    2  = @(45386.90608821759, 0.574341, "SylvAtri-b-m-5mn-haz-cos-44bdbaqj")
    3  = @(45386.90608967593, 0.5942499999999999, "SylvAtri-b-m-10mn-haz-cos-4mpnw3jt")
    4  = @(45386.90609125, 0.756752, "SylvAtri-ab-m-5mn-haz-cos-403ja14t")
    5  = @(45386.90609288195, 0.868096, "SylvAtri-ab-m-10mn-haz-cos-vgdwfgsx")
    6  = @(45386.90610482639, 0.549405, "TurdMeru-b-m-5mn-nex-cos-bfydabh4")
    7  = @(45386.90610868055, 0.419789, "TurdMeru-b-m-10mn-nex-cos-dqwkobyf")
    8  = @(45386.90611225695, 0.7795110000000001, "TurdMeru-ab-m-5mn-nex-cos-_s80wtoc")
    9  = @(45386.90612064815, 0.8188119999999999, "TurdMeru-ab-m-10mn-nex-cos-vunk_m1e")
    10 = @(45386.90613228009, 0.557665, "LuscMega-b-m-5mn-nex-cos-3othwfhn")
    11 = @(45386.90612702546, 0.358763, "LuscMega-b-m-10mn-nex-cos-oqdu6skh")
    12 = @(45386.90611974537, 0.680266, "LuscMega-ab-m-5mn-nex-cos-_ru9tilk")
    13 = @(45386.9061204051, 0.505812, "LuscMega-ab-m-10mn-nex-cos-drd26td0")
}

foreach ($row in $tuvData.Keys) {
    $values = $tuvData[$row]
    $ws.Range("T$row").Value = $values[0]
    $ws.Range("U$row").Value = $values[1]
    $ws.Range("V$row").Value = $values[2]
}

# ---------------------------------------------------------------------
# 2) "Synthesis" sheet : RunFolder column (AI) for rows 2..13
# ---------------------------------------------------------------------
$wsSynth = $wb.Worksheets.Item("Synthesis")

$aiData = @{
    2  = "SylvAtri-b-m-5mn-haz-cos-44bdbaqj"
    3  = "SylvAtri-b-m-10mn-haz-cos-4mpnw3jt"
    4  = "SylvAtri-ab-m-5mn-haz-cos-403ja14t"
    5  = "SylvAtri-ab-m-10mn-haz-cos-vgdwfgsx"
    6  = "TurdMeru-b-m-5mn-nex-cos-bfydabh4"
    7  = "TurdMeru-b-m-10mn-nex-cos-dqwkobyf"
    8  = "TurdMeru-ab-m-5mn-nex-cos-_s80wtoc"
    9  = "TurdMeru-ab-m-10mn-nex-cos-vunk_m1e"
    10 = "LuscMega-b-m-5mn-nex-cos-3othwfhn"
    11 = "LuscMega-b-m-10mn-nex-cos-oqdu6skh"
    12 = "LuscMega-ab-m-5mn-nex-cos-_ru9tilk"
    13 = "LuscMega-ab-m-10mn-nex-cos-drd26td0"
}

foreach ($row in $aiData.Keys) {
    $wsSynth.Range("AI$row").Value = $aiData[$row]
}

# ---------------------------------------------------------------------
# 3) "Computing platform" sheet : pyaudisam version and MCDS engine path
# ---------------------------------------------------------------------
$wsPlatform = $wb.Worksheets.Item("Computing platform")
$wsPlatform.Range("B10").Value = "1.0.3-b1"
$wsPlatform.Range("B11").Value = "C:/PortableApps/Distance 752/MCDS.exe"
